# Updates the cryptos price/volume table (Sheet1) to reflect the latest
# scrape, per commit "Updated cryptos list on Sat Mar  9 13:40:11 UTC 2024
# with GitHub Actions". Price values in column D are forced to text (via a
# leading apostrophe) so Excel doesn't coerce them to numbers and drop
# trailing zeros / treat multi-dot thousand separators as errors - matching
# the original inlineStr/text storage of that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'68.393.05"
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = "'3.921.85"
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'486.14"
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').Value = "'146.15"
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = "'0.167"
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('E11').Value = '  -2.58%  '
$ws.Range('D12').Value = "'43.16"
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = "'10.76"
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('D14').Value = "'4.541.16"
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').Value = "'3.934.16"
$ws.Range('E15').Value = '  -0.23%  '
$ws.Range('D16').Value = "'14.36"
$ws.Range('E16').Value = '  -5.39%  '
$ws.Range('E17').Value = '  -1.16%  '
$ws.Range('D18').Value = "'20.01"
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').Value = "'68.430.40"
$ws.Range('E20').Value = '  +1.11%  '
$ws.Range('D21').Value = "'433.56"
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').Value = "'15.27"
$ws.Range('E22').Value = '  +4.64%  '
$ws.Range('D23').Value = "'3.49"
$ws.Range('E23').Value = '  +1.99%  '
$ws.Range('D24').Value = "'88.36"
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').Value = "'11.50"
$ws.Range('E25').Value = '  +16.19%  '
$ws.Range('D26').Value = "'11.29"
$ws.Range('E26').Value = '  +11.45%  '
$ws.Range('D27').Value = "'3.66"
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = "'38.02"
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('D29').Value = "'5.70"
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').Value = "'714.20"
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('D31').Value = "'13.83"
$ws.Range('E31').Value = '  +2.92%  '
$ws.Range('E32').Value = '  -1.80%  '
$ws.Range('E33').Value = '  +4.00%  '
$ws.Range('D34').Value = "'6.15"
$ws.Range('E34').Value = '  +13.20%  '
$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = "'0.0₃0885"
$ws.Range('E35').Value = '  +4.05%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = "'41.35"
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('D37').Value = "'61.02"
$ws.Range('E37').Value = '  +3.28%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = "'0.146"
$ws.Range('E38').Value = '  -5.24%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = "'0.999"
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = "'3.02"
$ws.Range('E40').Value = '  +18.65%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = "'0.396"
$ws.Range('E41').Value = '  +16.23%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = "'0.0498"
$ws.Range('E42').Value = '  +4.34%  '
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('E44').Value = '  +4.57%  '
$ws.Range('D45').Value = "'3.39"
$ws.Range('E45').Value = '  +5.07%  '
$ws.Range('E46').Value = '  -1.32%  '
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('E48').Value = '  -1.92%  '
$ws.Range('E49').Value = '  -5.29%  '
$ws.Range('D50').Value = "'145.17"
$ws.Range('E50').Value = '  -2.60%  '
$ws.Range('D51').Value = "'0.0₆0339"
$ws.Range('E51').Value = '  +24.21%  '
